$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.145.51"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.492.45"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.95"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.83"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.38"
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.07"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "2.873.37"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "2.472.20"
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.92"
$ws.Range("E16").Value = "  -5.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "42.094.30"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").Value = "0.0₃0927"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.86"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.42"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.01"
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.76"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.92"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.86"
$ws.Range("E30").Value = "  -6.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.13"
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.70"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0765"
$ws.Range("E34").Value = "  -3.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.57"
$ws.Range("E35").Value = "  -8.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.06"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.37"
$ws.Range("E37").Value = "  -4.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.91"
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.09"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.24"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "2.009.19"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0289"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.75"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").Value = "2.714.33"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.78"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.184"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.35"
$ws.Range("E51").Value = "  -2.52%  "
